$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reshape the metrics table:
#   - columns C (SASA) / D (sum_SASA) are replaced by the monosaccharides /
#     motifs list columns that used to live in J / K
#   - column E (max_SASA) becomes "sasa" (value itself is unchanged)
#   - column F (flexibility) keeps its header but is recomputed
#   - column G (Q) becomes the boolean "has_multi_node_motifs" column
#   - columns H:L (theta, conformation, monosaccharides, motifs, class) are
#     dropped entirely
# ---------------------------------------------------------------------------

# Capture the values we still need (old monosaccharides / motifs columns)
# before we overwrite/delete anything.
$monosaccharides = @{}
$motifs = @{}
for ($r = 2; $r -le 5; $r++) {
    $monosaccharides[$r] = $ws.Cells.Item($r, 10).Value2   # column J
    $motifs[$r]          = $ws.Cells.Item($r, 11).Value2   # column K
}

# New header row
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
$ws.Range("F1").Value = "flexibility"
$ws.Range("G1").Value = "has_multi_node_motifs"

# New flexibility values (recomputed aggregation: max, sum)
$flexibility = @{
    2 = 0.5568992082145107
    3 = 0.7313890112283827
    4 = 0.1985667364886857
    5 = 0.1708058283471475
}

for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = $monosaccharides[$r]   # C: monosaccharides
    $ws.Cells.Item($r, 4).Value = $motifs[$r]            # D: motifs
    # E (sasa) already holds the correct max_SASA value, leave as-is
    $ws.Cells.Item($r, 6).Value = $flexibility[$r]       # F: flexibility
    $ws.Cells.Item($r, 7).Value = $false                 # G: has_multi_node_motifs
}

# Drop the now-unused columns H:L (theta, conformation, monosaccharides,
# motifs, class)
$ws.Range("H1:L5").EntireColumn.Delete()
